$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (A2:D12) by column A (time) in ascending order.
$dataRange = $ws.Range("A2:D12")
$keyRange = $ws.Range("A2:A12")

$dataRange.Sort($keyRange, 1)
